$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText {
    param($sheet, $cellRef, $text)
    $cell = $sheet.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value2 = $text
}

Set-CellText $ws "D2" "68.183.23"
Set-CellText $ws "E2" "  -0.24%  "
Set-CellText $ws "D3" "2.697.39"
Set-CellText $ws "E3" "  +1.83%  "
Set-CellText $ws "D4" "0.999"
Set-CellText $ws "E4" "  -0.28%  "
Set-CellText $ws "D5" "604.95"
Set-CellText $ws "E5" "  +1.52%  "
Set-CellText $ws "D6" "165.71"
Set-CellText $ws "E6" "  +4.55%  "
Set-CellText $ws "E7" "  +0.00%  "
Set-CellText $ws "D8" "0.551"
Set-CellText $ws "E8" "  +2.14%  "
Set-CellText $ws "D9" "2.696.54"
Set-CellText $ws "E9" "  +1.79%  "
Set-CellText $ws "E10" "  +3.43%  "
Set-CellText $ws "E11" "  +0.70%  "
Set-CellText $ws "D12" "0.361"
Set-CellText $ws "E12" "  +2.32%  "
Set-CellText $ws "D13" "5.27"
Set-CellText $ws "E13" "  +0.05%  "
Set-CellText $ws "D14" "28.10"
Set-CellText $ws "E14" "  +0.50%  "
Set-CellText $ws "D15" "3.189.05"
Set-CellText $ws "E15" "  +1.79%  "
Set-CellText $ws "D16" "0.0000186"
Set-CellText $ws "E16" "  +0.37%  "
Set-CellText $ws "D17" "68.181.42"
Set-CellText $ws "E17" "  -0.13%  "
Set-CellText $ws "D18" "2.692.86"
Set-CellText $ws "E18" "  +0.54%  "
Set-CellText $ws "E19" "  +0.86%  "
Set-CellText $ws "D20" "368.66"
Set-CellText $ws "E20" "  +1.48%  "
Set-CellText $ws "D21" "7.60"
Set-CellText $ws "E21" "  +1.44%  "
Set-CellText $ws "D22" "4.45"
Set-CellText $ws "E22" "  +0.63%  "
Set-CellText $ws "D23" "4.90"
Set-CellText $ws "E23" "  +2.49%  "
Set-CellText $ws "E24" "  -1.13%  "
Set-CellText $ws "D25" "72.36"
Set-CellText $ws "E25" "  -3.52%  "
Set-CellText $ws "E26" "  +0.10%  "
Set-CellText $ws "D27" "9.87"
Set-CellText $ws "E27" "  -0.61%  "
Set-CellText $ws "D28" "2.846.63"
Set-CellText $ws "E28" "  +2.12%  "
Set-CellText $ws "D29" "0.0000103"
Set-CellText $ws "E29" "  +1.43%  "
Set-CellText $ws "D30" "0.999"
Set-CellText $ws "E30" "  -0.07%  "
Set-CellText $ws "D31" "568.57"
Set-CellText $ws "E31" "  -1.47%  "
Set-CellText $ws "D32" "8.10"
Set-CellText $ws "E32" "  -0.36%  "
Set-CellText $ws "D33" "1.40"
Set-CellText $ws "E33" "  +0.47%  "
Set-CellText $ws "D34" "1.96"
Set-CellText $ws "E34" "  +4.55%  "
Set-CellText $ws "E35" "  +0.67%  "
Set-CellText $ws "D38" "19.78"
Set-CellText $ws "E38" "  +0.94%  "
Set-CellText $ws "D39" "157.73"
Set-CellText $ws "E39" "  -1.77%  "
Set-CellText $ws "D40" "0.376"
Set-CellText $ws "E40" "  +1.66%  "
Set-CellText $ws "D41" "5.36"
Set-CellText $ws "E41" "  +0.95%  "
Set-CellText $ws "D42" "1.85"
Set-CellText $ws "E42" "  -1.11%  "
Set-CellText $ws "D43" "17.97"
Set-CellText $ws "E43" "  +0.89%  "
Set-CellText $ws "D44" "2.56"
Set-CellText $ws "E44" "  -2.14%  "
Set-CellText $ws "D48" "0.593"
Set-CellText $ws "E48" "  +1.25%  "
Set-CellText $ws "D49" "155.00"
Set-CellText $ws "E49" "  -2.07%  "
Set-CellText $ws "D50" "3.88"
Set-CellText $ws "E50" "  +1.73%  "
Set-CellText $ws "D51" "1.75"
Set-CellText $ws "E51" "  +1.61%  "

# Row 36/37 swap: ImmutableX <-> FirstDigitalUSD
Set-CellText $ws "B36" "FirstDigitalUSD"
Set-CellText $ws "C36" "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-CellText $ws "D36" "0.999"
Set-CellText $ws "E36" "  -0.10%  "
Set-CellText $ws "B37" "ImmutableX"
Set-CellText $ws "C37" "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-CellText $ws "D37" "1.59"
Set-CellText $ws "E37" "  -2.20%  "

# Row 46/47 swap: BabyDogeCoin <-> OKB
Set-CellText $ws "B46" "OKB"
Set-CellText $ws "C46" "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-CellText $ws "D46" "40.65"
Set-CellText $ws "E46" "  +0.94%  "
Set-CellText $ws "B47" "BabyDogeCoin"
Set-CellText $ws "C47" "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-CellText $ws "D47" "0.0₆0306"
Set-CellText $ws "E47" "  -3.16%  "
